$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: highlight this run's inputs/result in yellow ---
$ws.Range("A15:K15").Interior.Color = 65535
$ws.Range("M15").Interior.Color = 65535

# --- Row 45 (new): another XGB run, same params as several others, lastvalue/d_mean, JCreeks ---
$ws.Range("A45").Value = -0.086
$ws.Range("B45").Value = 0.086
$ws.Range("C45").Value = 100
$ws.Range("D45").Value = 100
$ws.Range("E45").Value = 2
$ws.Range("F45").Value = 10
$ws.Range("G45").Value = 0.99
$ws.Range("H45").Value = 0.9
$ws.Range("I45").Value = 2
$ws.Range("K45").Value = 0.0144786
$ws.Range("L45").Value = "lastvalue, d_mean"
$ws.Range("M45").Value = "JCreeks"

# --- Row 47: tweak cut values, drop score, switch account ---
$ws.Range("A47").Value = -0.08575
$ws.Range("B47").Value = 0.08575
$ws.Range("K47").ClearContents()
$ws.Range("M47").Value = "jg6eb"

# --- Row 48: bump num_to_keep/quant, relabel trainer description ---
$ws.Range("F48").Value = 11
$ws.Range("G48").Value = 0.995
$ws.Range("L48").Value = "lastvalue, d_mean"

# --- Row 51 (new) ---
$ws.Range("A51").Value = -0.086
$ws.Range("B51").Value = 0.086
$ws.Range("C51").Value = 100
$ws.Range("D51").Value = 100
$ws.Range("E51").Value = 2
$ws.Range("F51").Value = 10
$ws.Range("G51").Value = 0.995
$ws.Range("H51").Value = 0.9
$ws.Range("I51").Value = 2
$ws.Range("K51").Value = 0.014439
$ws.Range("L51").Value = "lastvalue, d_mean, mean residual"
$ws.Range("M51").Value = "JC113"

# --- Row 52 (new) ---
$ws.Range("A52").Value = -0.08575
$ws.Range("B52").Value = 0.08575
$ws.Range("C52").Value = 100
$ws.Range("D52").Value = 100
$ws.Range("E52").Value = 2
$ws.Range("F52").Value = 10
$ws.Range("G52").Value = 0.995
$ws.Range("H52").Value = 0.9
$ws.Range("I52").Value = 2
$ws.Range("L52").Value = "lastvalue, d_mean, mean residual"
$ws.Range("M52").Value = "JC113"

# --- Row 53 (new) ---
$ws.Range("A53").Value = -0.086
$ws.Range("B53").Value = 0.086
$ws.Range("C53").Value = 100
$ws.Range("D53").Value = 100
$ws.Range("E53").Value = 2
$ws.Range("F53").Value = 11
$ws.Range("G53").Value = 0.995
$ws.Range("H53").Value = 0.9
$ws.Range("I53").Value = 2
$ws.Range("K53").Value = "OT"
$ws.Range("L53").Value = "lastvalue, d_mean, mean residual"
$ws.Range("M53").Value = "JC113"

# --- Row 54 (new) ---
$ws.Range("A54").Value = -0.086
$ws.Range("B54").Value = 0.086
$ws.Range("C54").Value = 100
$ws.Range("D54").Value = 100
$ws.Range("E54").Value = 2
$ws.Range("F54").Value = 10
$ws.Range("G54").Value = 0.99
$ws.Range("H54").Value = 0.9
$ws.Range("I54").Value = 2
$ws.Range("L54").Value = "lastvalue, d_mean, mean residual"
$ws.Range("M54").Value = "jg6ebuva"

# --- Row 56 (new, row 55 left blank) ---
$ws.Range("A56").Value = -0.086
$ws.Range("B56").Value = 0.086
$ws.Range("C56").Value = 100
$ws.Range("D56").Value = 100
$ws.Range("E56").Value = 2
$ws.Range("F56").Value = 10
$ws.Range("G56").Value = 0.99
$ws.Range("H56").Value = 0.9
$ws.Range("I56").Value = 2
$ws.Range("K56").Value = -0.0191153
$ws.Range("L56").Value = "lastvalue, d_mean, mean residual, trainer=residual, L3 fit train.y XGB ensemble"
$ws.Range("M56").Value = "jg6eb"

# --- Row 59 (new, rows 57-58 left blank) ---
$ws.Range("A59").Value = -0.086
$ws.Range("B59").Value = 0.086
$ws.Range("C59").Value = 100
$ws.Range("D59").Value = 100
$ws.Range("E59").Value = 2
$ws.Range("F59").Value = 10
$ws.Range("G59").Value = 0.99
$ws.Range("H59").Value = 0.9
$ws.Range("I59").Value = 2
$ws.Range("K59").Value = 0.0115586
$ws.Range("L59").Value = "lastvalue, d_mean, mean residual, trainer=residual, L3 fit train.y rft stacking"
$ws.Range("M59").Value = "jg6ebuva"

# --- Move selection to where editing finished ---
$ws.Range("L48").Select()
